# "Erase first row instead of iteration on trials"
# Adds a new "sub_num" / "subject number" column (AU) to the header rows,
# matching the style of the existing header cells, and moves the sheet's
# view/selection to the newly added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string header entries for column AU.
$ws.Range("AU1").Value = "sub_num"
$ws.Range("AU2").Value = "subject number"

# Match the wrap-text style used by every other header cell (style index 1).
$ws.Range("AU1:AU2").WrapText = $true

# Scroll the view toward the new column and select the cell below it.
$excel.ActiveWindow.ScrollColumn = 35
$ws.Range("AU3").Select()
